# "add workbook buffer reader"
# Populate Sheet2 with a second batch of card records (mirroring the
# layout already present on Sheet1), then leave the workbook positioned
# with Sheet2 active/selected, as a user would after entering this data.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Header row (identical to Sheet1's header row)
$ws2.Range("A1").Value = "cardNo"
$ws2.Range("B1").Value = "userName"
$ws2.Range("C1").Value = "cardPwd"
$ws2.Range("D1").Value = "cardType.id"
$ws2.Range("E1").Value = "cardType.cardBean.cardNo"

# Row 2
$ws2.Range("A2").Value = "2-111111"
$ws2.Range("B2").Value = "马日胜"
$ws2.Range("C2").Value = 111111111
$ws2.Range("D2").Value = 5000
$ws2.Range("E2").Value = 6666

# Row 3
$ws2.Range("A3").Value = "2-222222"
$ws2.Range("B3").Value = "马胜日"
$ws2.Range("C3").Value = 222222
$ws2.Range("D3").Value = 5000
$ws2.Range("E3").Value = 6677

# Row 4
$ws2.Range("A4").Value = "2-33333"
$ws2.Range("B4").Value = "日胜马"
$ws2.Range("C4").Value = 3333333
$ws2.Range("D4").Value = 5000
$ws2.Range("E4").Value = 77

# Row 5
$ws2.Range("A5").Value = "2-44444"
$ws2.Range("B5").Value = "日马胜"
$ws2.Range("C5").Value = 4444444
$ws2.Range("D5").Value = 5000
$ws2.Range("E5").Value = 88

# Sheet1 ends up with its used rows (1:5) selected instead of the single
# cell G3 that was previously selected there.
$ws1.Activate()
$ws1.Rows("1:5").Select()

# Sheet2 becomes the active sheet/tab, with G10 selected.
$ws2.Activate()
$ws2.Range("G10").Select()
